$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-31 Wednesday" "2024-08-01 Thursday"

Replace-Text "665×6=3990" "743×5=3715"
Replace-Text "205×4=820" "541×7=3787"
Replace-Text "182×8=1456" "944×7=6608"
Replace-Text "244×8=1952" "892×8=7136"
Replace-Text "861×2=1722" "190×2=380"

Replace-Text "400×4=1600" "721×9=6489"
Replace-Text "883×3=2649" "843×9=7587"
Replace-Text "520×6=3120" "331×7=2317"
Replace-Text "412×9=3708" "365×2=730"
Replace-Text "227×3=681" "127×2=254"

Replace-Text "380×5=1900" "982×2=1964"
Replace-Text "432×2=864" "535×3=1605"
Replace-Text "750×9=6750" "616×4=2464"
Replace-Text "854×7=5978" "612×2=1224"
Replace-Text "517×7=3619" "635×3=1905"

Replace-Text "517×5=2585" "847×4=3388"
Replace-Text "603×8=4824" "388×6=2328"
Replace-Text "179×3=537" "848×5=4240"
Replace-Text "361×7=2527" "251×5=1255"
Replace-Text "789×2=1578" "930×9=8370"

Replace-Text "250×3=750" "701×6=4206"
Replace-Text "536×7=3752" "497×7=3479"
Replace-Text "948×6=5688" "262×7=1834"
Replace-Text "603×2=1206" "500×7=3500"
Replace-Text "792×3=2376" "834×3=2502"
